$wb = $excel.ActiveWorkbook

# --- Sheet "右1-..." (4th sheet in tab order) gets new column C and a new name/layout ---
$ws4 = $wb.Worksheets.Item(4)

# Rewrite the header row: A1 becomes blank, B1/C1 get the two new headers.
$ws4.Range("A1").ClearContents()
$ws4.Range("B1").Value = "海水养殖（万吨)`n"
$ws4.Range("C1").Value = "淡水养殖(万吨)`n"

# Rewrite the data rows (A stays as the year; B/C get the new sea/fresh water figures).
$ws4.Range("B2").Value = 29.42
$ws4.Range("C2").Value = 36.69
$ws4.Range("B3").Value = 27.1
$ws4.Range("C3").Value = 35.44
$ws4.Range("B4").Value = 26.65
$ws4.Range("C4").Value = 35.09
$ws4.Range("B5").Value = 27.3
$ws4.Range("C5").Value = 35.18
$ws4.Range("B6").Value = 26.26
$ws4.Range("C6").Value = 41.08

# Re-apply the existing "year column" style (already used by A2:A6) across the whole
# A1:C6 block so the new header/column match the surrounding formatting.
$ws4.Range("A2").Copy()
$ws4.Range("A1:C6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws4.Rows.Item(1).RowHeight = 39.75

# Rename the sheet to reflect the new chart contents.
$ws4.Name = "右1-近年养殖水产品产量及构成"

# Make this sheet the active one (was "左2-...", sheet index 2).
$ws4.Activate()
$ws4.Range("H8").Select()
